$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Switch workbook calculation mode from manual back to automatic
$wb.Application.Calculation = -4105  # xlCalculationAutomatic

# Row 22: new "egreso" record
$ws.Range("A22").Value = "2017-11-07-20:00:00"
$ws.Range("B22").Value = 32314054
$ws.Range("C22").Value = "Perez"
$ws.Range("D22").Value = "Pedro"
$ws.Range("E22").Value = "HTW327"
$ws.Range("F22").Value = "S/D"
$ws.Range("G22").Value = "S/D"
$ws.Range("H22").Value = 2
$ws.Range("I22").Value = 121
$ws.Range("J22").Value = "DNI"

# Row 23: new "egreso" record
$ws.Range("A23").Value = "2017-11-07-20:00:10"
$ws.Range("B23").Value = 25745226
$ws.Range("C23").Value = "Hernandez"
$ws.Range("D23").Value = "Cristian"
$ws.Range("E23").Value = "HVS 839"
$ws.Range("F23").Value = "S/D"
$ws.Range("G23").Value = "S/D"
$ws.Range("H23").Value = 3
$ws.Range("I23").Value = 122
$ws.Range("J23").Value = "DNI"

# Move selection to match the post-edit active cell
$ws.Range("H22").Select()
